# Scheduled-runner update: refresh cached market-board pricing figures
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for a
# batch of Leve rows across the ALC / ARM / BSM / CRP / CUL / GSM / WVR
# sheets, per the latest price-data pull.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3918
$ws.Range("J32").Value = 3918
$ws.Range("L32").Value = 3918
$ws.Range("N32").Value = -4570

$ws.Range("H64").Value = 4517
$ws.Range("I64").Value = 4073
$ws.Range("J64").Value = 5072
$ws.Range("K64").Value = 4073
$ws.Range("L64").Value = 5072
$ws.Range("M64").Value = -3825
$ws.Range("N64").Value = -5568

$ws.Range("H67").Value = 4517
$ws.Range("I67").Value = 4073
$ws.Range("J67").Value = 5072
$ws.Range("K67").Value = 4073
$ws.Range("L67").Value = 5072
$ws.Range("M67").Value = -3215
$ws.Range("N67").Value = -6788

$ws.Range("H76").Value = 6000
$ws.Range("J76").Value = 6000
$ws.Range("L76").Value = 6000
$ws.Range("N76").Value = -6630

$ws.Range("H79").Value = 6000
$ws.Range("J79").Value = 6000
$ws.Range("L79").Value = 6000
$ws.Range("N79").Value = -8184

$ws.Range("H113").Value = 4116.6665
$ws.Range("J113").Value = 4250
$ws.Range("L113").Value = 4250
$ws.Range("N113").Value = -10758

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5637
$ws.Range("I45").Value = 6716.1665
$ws.Range("J45").Value = 2399.5
$ws.Range("K45").Value = 6716.1665
$ws.Range("L45").Value = 2399.5
$ws.Range("M45").Value = -6339.1665
$ws.Range("N45").Value = -3153.5

$ws.Range("H132").Value = 1231.1666
$ws.Range("I132").Value = 1231.1666
$ws.Range("K132").Value = 3693.4998
$ws.Range("M132").Value = -1163.4998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1546.8
$ws.Range("I20").Value = 1445
$ws.Range("J20").Value = 1699.5
$ws.Range("K20").Value = 1445
$ws.Range("L20").Value = 1699.5
$ws.Range("M20").Value = -1198
$ws.Range("N20").Value = -2193.5

$ws.Range("H86").Value = 2405
$ws.Range("I86").Value = 2031.25
$ws.Range("J86").Value = 3900
$ws.Range("K86").Value = 2031.25
$ws.Range("L86").Value = 3900
$ws.Range("M86").Value = -908.25
$ws.Range("N86").Value = -6146

$ws.Range("H89").Value = 2405
$ws.Range("I89").Value = 2031.25
$ws.Range("J89").Value = 3900
$ws.Range("K89").Value = 10156.25
$ws.Range("L89").Value = 19500
$ws.Range("M89").Value = -4540.25
$ws.Range("N89").Value = -30732

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 575.5714
$ws.Range("I10").Value = 355.25
$ws.Range("J10").Value = 869.3333
$ws.Range("K10").Value = 355.25
$ws.Range("L10").Value = 869.3333
$ws.Range("M10").Value = -216.25
$ws.Range("N10").Value = -1147.3333

$ws.Range("H19").Value = 63.090908
$ws.Range("I19").Value = 63.090908
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 63.090908
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = 106.909092
$ws.Range("M19").ClearContents()

$ws.Range("H24").Value = 63.090908
$ws.Range("I24").Value = 63.090908
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 63.090908
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = 106.909092
$ws.Range("M24").ClearContents()

$ws.Range("H62").Value = 3872.25
$ws.Range("I62").Value = 2329.6667
$ws.Range("K62").Value = 2329.6667
$ws.Range("M62").Value = -1705.6667

$ws.Range("H65").Value = 3872.25
$ws.Range("I65").Value = 2329.6667
$ws.Range("K65").Value = 11648.3335
$ws.Range("M65").Value = -8528.333500000001

$ws.Range("H99").Value = 3333
$ws.Range("I99").Value = 999.5
$ws.Range("K99").Value = 999.5
$ws.Range("M99").Value = 498.5

$ws.Range("H126").Value = 3333
$ws.Range("I126").Value = 999.5
$ws.Range("K126").Value = 2998.5
$ws.Range("M126").Value = -528.5

$ws.Range("H132").Value = 1703.7142
$ws.Range("I132").Value = 1622.12
$ws.Range("K132").Value = 4866.36
$ws.Range("M132").Value = -2336.36

$ws.Range("H134").Value = 1051.8462
$ws.Range("I134").Value = 1087.3636
$ws.Range("K134").Value = 3262.0908
$ws.Range("M134").Value = -727.0907999999999

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59954
$ws.Range("J37").Value = 59954
$ws.Range("L37").Value = 179862
$ws.Range("N37").Value = -180086

$ws.Range("H39").Value = 9420.75
$ws.Range("J39").Value = 9894.666999999999
$ws.Range("L39").Value = 29684.001
$ws.Range("N39").Value = -30272.001

$ws.Range("H113").Value = 1497.8948
$ws.Range("I113").Value = 1382.5
$ws.Range("J113").Value = 1511.4706
$ws.Range("K113").Value = 4147.5
$ws.Range("L113").Value = 4534.4118
$ws.Range("M113").Value = -1977.5
$ws.Range("N113").Value = -8874.4118

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("N122").Value = 0
$ws.Range("L122").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4979.4
$ws.Range("I70").Value = 4425
$ws.Range("J70").Value = 6088.2
$ws.Range("K70").Value = 4425
$ws.Range("L70").Value = 6088.2
$ws.Range("M70").Value = -4155
$ws.Range("N70").Value = -6628.2

$ws.Range("H73").Value = 4979.4
$ws.Range("I73").Value = 4425
$ws.Range("J73").Value = 6088.2
$ws.Range("K73").Value = 4425
$ws.Range("L73").Value = 6088.2
$ws.Range("M73").Value = -3489
$ws.Range("N73").Value = -7960.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 18850000
$ws.Range("I2").Value = 31000000
$ws.Range("J2").Value = 10750000
$ws.Range("K2").Value = 31000000
$ws.Range("L2").Value = 10750000
$ws.Range("M2").Value = -30999888
$ws.Range("N2").Value = -10750224

$ws.Range("H3").Value = 4999.5
$ws.Range("J3").Value = 4999.5
$ws.Range("L3").Value = 4999.5
$ws.Range("N3").Value = -5227.5

$ws.Range("H4").Value = 673518.5
$ws.Range("I4").Value = 1333703.6
$ws.Range("J4").Value = 13333.333
$ws.Range("K4").Value = 1333703.6
$ws.Range("L4").Value = 13333.333
$ws.Range("M4").Value = -1333590.6
$ws.Range("N4").Value = -13559.333

$ws.Range("H5").Value = 20999.4
$ws.Range("J5").Value = 20999.4
$ws.Range("L5").Value = 20999.4
$ws.Range("N5").Value = -21223.4

$ws.Range("H7").Value = 9500
$ws.Range("I7").Value = 1250
$ws.Range("J7").Value = 17750
$ws.Range("K7").Value = 1250
$ws.Range("L7").Value = 17750
$ws.Range("M7").Value = -1137
$ws.Range("N7").Value = -17976

$ws.Range("H9").Value = 1100
$ws.Range("J9").Value = 1100
$ws.Range("L9").Value = 1100
$ws.Range("N9").Value = -1380

$ws.Range("H10").Value = 643333.3
$ws.Range("J10").Value = 465000
$ws.Range("L10").Value = 465000
$ws.Range("N10").Value = -465338

$ws.Range("H14").Value = 1500
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 2000
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 2000
$ws.Range("M14").Value = -832
$ws.Range("N14").Value = -2336

$ws.Range("H100").Value = 5556340.5
$ws.Range("I100").Value = 6250861.5
$ws.Range("J100").Value = 174
$ws.Range("K100").Value = 12501723
$ws.Range("L100").Value = 348
$ws.Range("M100").Value = -12501182
$ws.Range("N100").Value = -1430

$ws.Range("H116").Value = 42840
$ws.Range("J116").Value = 42840
$ws.Range("L116").Value = 42840
$ws.Range("N116").Value = -52018

$ws.Range("H117").Value = 96759.664
$ws.Range("J117").Value = 96759.664
$ws.Range("L117").Value = 96759.664
$ws.Range("N117").Value = -105937.664

$ws.Range("H132").Value = 2461
$ws.Range("I132").Value = 2559.7144
$ws.Range("K132").Value = 7679.1432
$ws.Range("M132").Value = -5149.1432
